$wb = $excel.ActiveWorkbook

# --- Device sheet: remove the "100, Floating Bus" section header row ---
# In the original workbook this was row 16 (an empty A16 + B16 = "100, Floating Bus"
# header, mirroring the "0-10, Synchronous Generator" / "10-19, Grid-Following VSI" /
# "90-99, Passive Load" section headers above it). Deleting the whole row shifts
# every following row up by one and drops the now-unused shared string automatically.
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Rows("16:16").Delete()

# --- Update the saved selections on each sheet to match the new view state ---
# NetworkLine: selection moved to D20 (sheet stays non-active in the final state)
$wsNetworkLine = $wb.Worksheets.Item("NetworkLine")
$wsNetworkLine.Range("D20").Select()

# Device: selection moved to B23 (sheet stays non-active in the final state)
$wsDevice.Range("B23").Select()

# PowerFlow: selection moved to F16, and this becomes the active/selected tab
$wsPowerFlow = $wb.Worksheets.Item("PowerFlow")
$wsPowerFlow.Activate()
$wsPowerFlow.Range("F16").Select()
